$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry per-record data subject to the weekly reshuffle.
$cols = @("D","H","J","K","L","M","N","O","P","Q")

# Rows that participate in the reorder (5 and 22 are left untouched).
$rows = @(2,3,4,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,23,24,25)

# Snapshot the original values for every (row, col) before any writes happen,
# since the reorder is a permutation and later writes must not read data
# that has already been overwritten.
$snapshot = @{}
foreach ($r in $rows) {
    $rec = @{}
    foreach ($c in $cols) {
        $rec[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rec
}

# Destination row -> source row (source row's original record is copied into
# destination row).
$mapping = @{
    2  = 7
    3  = 8
    4  = 19
    6  = 16
    7  = 20
    8  = 3
    9  = 14
    10 = 15
    11 = 21
    12 = 25
    13 = 23
    14 = 13
    15 = 9
    16 = 4
    17 = 6
    18 = 2
    19 = 24
    20 = 17
    21 = 10
    23 = 18
    24 = 11
    25 = 12
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $srcRec = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcRec[$c]
    }
}
